$d = $word.ActiveDocument

# The document contains three "<id>...</id>" markers that were each split
# across three separate runs (one for "<id>", one for the bare identifier
# text, and one for "</id>"). They should be collapsed into a single run
# per marker. Two of the three keep their identifier text unchanged while
# the third's identifier text changes from "p042r_a3" to "p042r_3".

$d.Content.Find.Execute("<id>p042r_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p042r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p042r_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p042r_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p042r_a3</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p042r_3</id>", 2) | Out-Null
